$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so numeric-looking
# strings (e.g. "1.00", "0.0330") are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '61.112.44'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = '2.884.42'
$ws.Range('E3').Value = '  -1.39%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '589.61'
$ws.Range('E5').Value = '  -0.26%  '
$ws.Range('D6').Value = '139.41'
$ws.Range('E6').Value = '  -5.37%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = '2.883.32'
$ws.Range('E8').Value = '  -1.44%  '
$ws.Range('D9').Value = '0.493'
$ws.Range('E9').Value = '  -2.81%  '
$ws.Range('D10').Value = '7.12'
$ws.Range('E10').Value = '  +2.39%  '
$ws.Range('D11').Value = '0.138'
$ws.Range('E11').Value = '  -4.09%  '
$ws.Range('D12').Value = '0.429'
$ws.Range('E12').Value = '  -2.70%  '
$ws.Range('D13').Value = '0.0000218'
$ws.Range('E13').Value = '  -3.69%  '
$ws.Range('D14').Value = '32.24'
$ws.Range('E14').Value = '  -4.33%  '
$ws.Range('E15').Value = '  -0.43%  '
$ws.Range('D16').Value = '3.360.48'
$ws.Range('E16').Value = '  -1.45%  '
$ws.Range('D17').Value = '60.963.75'
$ws.Range('E17').Value = '  +0.01%  '
$ws.Range('D18').Value = '2.885.78'
$ws.Range('E18').Value = '  -1.37%  '
$ws.Range('D19').Value = '6.51'
$ws.Range('E19').Value = '  -2.92%  '
$ws.Range('D20').Value = '426.52'
$ws.Range('E20').Value = '  -1.44%  '
$ws.Range('D21').Value = '13.14'
$ws.Range('E21').Value = '  -2.15%  '
$ws.Range('D22').Value = '0.655'
$ws.Range('E22').Value = '  -3.50%  '
$ws.Range('D23').Value = '6.90'
$ws.Range('E23').Value = '  -2.97%  '
$ws.Range('D24').Value = '79.95'
$ws.Range('E24').Value = '  -1.78%  '
$ws.Range('D25').Value = '10.49'
$ws.Range('E25').Value = '  -3.53%  '
$ws.Range('B26').Value = 'Fetch.AI'
$ws.Range('C26').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D26').Value = '2.07'
$ws.Range('E26').Value = '  -6.29%  '
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('D28').Value = '11.44'
$ws.Range('E28').Value = '  -4.02%  '
$ws.Range('D29').Value = '2.55'
$ws.Range('E29').Value = '  -2.48%  '
$ws.Range('D30').Value = '2.07'
$ws.Range('E30').Value = '  -8.81%  '
$ws.Range('D31').Value = '6.63'
$ws.Range('E31').Value = '  -5.31%  '
$ws.Range('E32').Value = '  -0.08%  '
$ws.Range('D33').Value = '25.65'
$ws.Range('E33').Value = '  -3.85%  '
$ws.Range('E34').Value = '  -5.12%  '
$ws.Range('D35').Value = '0.0₃0842'
$ws.Range('E35').Value = '  -2.60%  '
$ws.Range('E36').Value = '  -4.56%  '
$ws.Range('D37').Value = '5.43'
$ws.Range('E37').Value = '  -3.80%  '
$ws.Range('D38').Value = '48.96'
$ws.Range('E38').Value = '  -2.16%  '
$ws.Range('D39').Value = '2.80'
$ws.Range('E39').Value = '  -7.18%  '
$ws.Range('D40').Value = '1.90'
$ws.Range('E40').Value = '  -4.13%  '
$ws.Range('D41').Value = '8.33'
$ws.Range('E41').Value = '  -2.58%  '
$ws.Range('E42').Value = '  -5.20%  '
$ws.Range('D43').Value = '0.267'
$ws.Range('E43').Value = '  -5.88%  '
$ws.Range('D44').Value = '38.61'
$ws.Range('E44').Value = '  -6.60%  '
$ws.Range('D45').Value = '2.662.02'
$ws.Range('E45').Value = '  -1.67%  '
$ws.Range('D46').Value = '131.49'
$ws.Range('E46').Value = '  -1.52%  '
$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D47').Value = '357.05'
$ws.Range('E47').Value = '  -5.26%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').Value = '0.0330'
$ws.Range('E48').Value = '  -4.39%  '
$ws.Range('E49').Value = '  +0.08%  '
$ws.Range('E50').Value = '  -4.02%  '
$ws.Range('D51').Value = '22.30'
$ws.Range('E51').Value = '  -6.93%  '
